# Auto-generated edit script: updates cryptos list values per commit
# "Updated cryptos list on Tue Feb 20 04:31:56 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '51.753.72'
$ws.Cells.Item(2, 5).Value = '  -0.83%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.918.85'
$ws.Cells.Item(3, 5).Value = '  +0.88%  '
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '353.43'
$ws.Cells.Item(5, 5).Value = '  +0.23%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '109.54'
$ws.Cells.Item(6, 5).Value = '  -2.25%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.572'
$ws.Cells.Item(7, 5).Value = '  +1.60%  '
$ws.Cells.Item(8, 5).Value = '  +0.00%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.626'
$ws.Cells.Item(9, 5).Value = '  +0.50%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '39.13'
$ws.Cells.Item(10, 5).Value = '  -2.38%  '
$ws.Cells.Item(11, 5).Value = '  +2.88%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '19.64'
$ws.Cells.Item(13, 5).Value = '  -1.77%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.85'
$ws.Cells.Item(14, 5).Value = '  +0.78%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '3.373.31'
$ws.Cells.Item(15, 5).Value = '  +1.02%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '2.918.02'
$ws.Cells.Item(16, 5).Value = '  +1.87%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.974'
$ws.Cells.Item(17, 5).Value = '  -1.72%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '51.756.56'
$ws.Cells.Item(18, 5).Value = '  -0.80%  '
$ws.Cells.Item(19, 2).Value = 'ImmutableX'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.28'
$ws.Cells.Item(19, 5).Value = '  -2.13%  '
$ws.Cells.Item(20, 2).Value = 'Uniswap'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.50'
$ws.Cells.Item(20, 5).Value = '  -3.04%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.86'
$ws.Cells.Item(21, 5).Value = '  -4.01%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.0₃0979'
$ws.Cells.Item(22, 5).Value = '  -0.28%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '70.63'
$ws.Cells.Item(23, 5).Value = '  -0.56%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '269.57'
$ws.Cells.Item(24, 5).Value = '  -0.24%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.79'
$ws.Cells.Item(25, 5).Value = '  -0.01%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.183'
$ws.Cells.Item(26, 5).Value = '  +11.70%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '27.03'
$ws.Cells.Item(27, 5).Value = '  +2.16%  '
$ws.Cells.Item(28, 5).Value = '  +0.05%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.29'
$ws.Cells.Item(29, 5).Value = '  +12.84%  '
$ws.Cells.Item(30, 5).Value = '  +10.89%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '10.50'
$ws.Cells.Item(31, 5).Value = '  -0.05%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '38.53'
$ws.Cells.Item(32, 5).Value = '  -0.92%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.00'
$ws.Cells.Item(33, 5).Value = '  -2.05%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '51.98'
$ws.Cells.Item(34, 5).Value = '  -2.56%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0437'
$ws.Cells.Item(35, 5).Value = '  -4.86%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.998'
$ws.Cells.Item(36, 5).Value = '  +0.05%  '
$ws.Cells.Item(37, 5).Value = '  -16.13%  '
$ws.Cells.Item(38, 5).Value = '  -3.20%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '18.31'
$ws.Cells.Item(39, 5).Value = '  -1.77%  '
$ws.Cells.Item(40, 5).Value = '  -2.18%  '
$ws.Cells.Item(41, 5).Value = '  +4.24%  '
$ws.Cells.Item(42, 5).Value = '  +2.46%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '22.66'
$ws.Cells.Item(43, 5).Value = '  -0.52%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '121.02'
$ws.Cells.Item(44, 5).Value = '  -0.63%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.17'
$ws.Cells.Item(45, 5).Value = '  -1.60%  '
$ws.Cells.Item(46, 5).Value = '  +0.83%  '
$ws.Cells.Item(47, 5).Value = '  -4.34%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.133.57'
$ws.Cells.Item(48, 5).Value = '  -3.28%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.248'
$ws.Cells.Item(49, 5).Value = '  -7.73%  '
$ws.Cells.Item(50, 5).Value = '  +2.63%  '
$ws.Cells.Item(51, 2).Value = 'FraxShare'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '9.03'
$ws.Cells.Item(51, 5).Value = '  -1.74%  '
